$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 426, shifting existing rows 426-485 down to 429-488.
$ws.Range("A426:T428").EntireRow.Insert()

# Populate the 3 newly-inserted rows with new weekly price records
# (same market/product/variety metadata as the block that used to start at row 426,
#  now pushed down to row 429, but with a new date and updated prices).

# Row 426 - Especial
$ws.Cells.Item(426, 1).Value = 4
$ws.Cells.Item(426, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(426, 3).Value = "Los Lagos"
$ws.Cells.Item(426, 4).Value = 45034
$ws.Cells.Item(426, 5).Value = 10
$ws.Cells.Item(426, 6).Value = "Fruta"
$ws.Cells.Item(426, 7).Value = 100101
$ws.Cells.Item(426, 8).Value = "Berries"
$ws.Cells.Item(426, 9).Value = 100101007
$ws.Cells.Item(426, 10).Value = "Kiwi"
$ws.Cells.Item(426, 11).Value = "Hayward"
$ws.Cells.Item(426, 12).Value = "Especial"
$ws.Cells.Item(426, 13).Value = 300
$ws.Cells.Item(426, 14).Value = 21000
$ws.Cells.Item(426, 15).Value = 21000
$ws.Cells.Item(426, 16).Value = 21000
$ws.Cells.Item(426, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(426, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(426, 19).Value = 1400
$ws.Cells.Item(426, 20).Value = 15

# Row 427 - Primera
$ws.Cells.Item(427, 1).Value = 4
$ws.Cells.Item(427, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(427, 3).Value = "Los Lagos"
$ws.Cells.Item(427, 4).Value = 45034
$ws.Cells.Item(427, 5).Value = 10
$ws.Cells.Item(427, 6).Value = "Fruta"
$ws.Cells.Item(427, 7).Value = 100101
$ws.Cells.Item(427, 8).Value = "Berries"
$ws.Cells.Item(427, 9).Value = 100101007
$ws.Cells.Item(427, 10).Value = "Kiwi"
$ws.Cells.Item(427, 11).Value = "Hayward"
$ws.Cells.Item(427, 12).Value = "Primera"
$ws.Cells.Item(427, 13).Value = 300
$ws.Cells.Item(427, 14).Value = 19000
$ws.Cells.Item(427, 15).Value = 19000
$ws.Cells.Item(427, 16).Value = 19000
$ws.Cells.Item(427, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(427, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(427, 19).Value = 1267
$ws.Cells.Item(427, 20).Value = 15

# Row 428 - Segunda
$ws.Cells.Item(428, 1).Value = 4
$ws.Cells.Item(428, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(428, 3).Value = "Los Lagos"
$ws.Cells.Item(428, 4).Value = 45034
$ws.Cells.Item(428, 5).Value = 10
$ws.Cells.Item(428, 6).Value = "Fruta"
$ws.Cells.Item(428, 7).Value = 100101
$ws.Cells.Item(428, 8).Value = "Berries"
$ws.Cells.Item(428, 9).Value = 100101007
$ws.Cells.Item(428, 10).Value = "Kiwi"
$ws.Cells.Item(428, 11).Value = "Hayward"
$ws.Cells.Item(428, 12).Value = "Segunda"
$ws.Cells.Item(428, 13).Value = 300
$ws.Cells.Item(428, 14).Value = 15000
$ws.Cells.Item(428, 15).Value = 15000
$ws.Cells.Item(428, 16).Value = 15000
$ws.Cells.Item(428, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(428, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(428, 19).Value = 1000
$ws.Cells.Item(428, 20).Value = 15
